$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1329217","Architectural Design Assistant","Београд, Србија","No","0 applicants","9 - 12 Weeks","Abovus architects"),
    @("1329081","Accelerate Romania - Web Developer","Cluj-Napoca, Romania","No","7 applicants","9 - 12 Weeks","GIS Solutions"),
    @("1329079","Accelerate Romania - YOUNG CIVIL & GEOTECHNICAL ENGINEER","Cluj-Napoca, Romania","No","2 applicants","9 - 12 Weeks","Expert Proiect"),
    @("1327539","Graphic Designer","União das freguesias de Cascais e Estoril, Portugal","No","42 applicants","3 - 6 Months","Dark Cloud"),
    @("1325528","Accelerate Romania | Spanish-Speaking Purchasing & Negotiation Specialist (EU Citizenship Required)","Iași, Romania","No","26 applicants","6 - 18 Months","Veo Wordwide Services - Iași"),
    @("1325527","Accelerate Romania | Spanish-Speaking Purchasing & Negotiation Specialist (EU Citizenship Required)","Brașov, Romania","No","20 applicants","6 - 18 Months","Veo Wordwide Services - Iași"),
    @("1325524","Accelerate Romania | Turkish-Speaking Purchasing & Negotiation Specialist (EU Citizenship Required)","Constanța, Romania","No","35 applicants","6 - 18 Months","Veo Wordwide Services - Iași"),
    @("1317005","Sales","Manisa, Yunusemre/Manisa, Türkiye","No","105 applicants","9 - 12 Weeks","TOYO MATBAA MÜREKKEPLERİ SANAYİ VE TİCARET ANONİM ŞİRKETİ"),
    @("1314780","Data Delivery & Analysis Intern (EU Preferred)","Heerlen, Nederland","No","217 applicants","6 - 18 Months","APG Heerlen"),
    @("1307425","Service executive II","Naucalpan de Juárez, Mexico","No","32 applicants","6 - 18 Months","Segmenta S.C.")
)

$row = 2
foreach ($item in $data) {
    $id = $item[0]
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $id
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $ws.Cells.Item($row, 5).Value = $item[3]
    $ws.Cells.Item($row, 6).Value = $item[4]
    $ws.Cells.Item($row, 7).Value = $item[5]
    $ws.Cells.Item($row, 8).Value = $item[6]
    $row = $row + 1
}

$ws.Range("E5").ClearFormats()
$ws.Range("E6").ClearFormats()

$ws.Range("C1").ColumnWidth = 101.17
$ws.Range("D1").ColumnWidth = 53.17
$ws.Range("H1").ColumnWidth = 59.17
